$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "MS_DEF"
